$d = $word.ActiveDocument

# The first paragraph in the document holds the AFFARS placeholder ID text:
#   "**ID__AFFARS_pgi_5319_topic_5__ID** " (placeholder run + trailing-space run)
$p1 = $d.Paragraphs.Item(1)

# 1. Add a paragraph border (top/left/bottom/right) with 5-twip spacing,
#    matching the pBdr already present on the third paragraph.
$borders = $p1.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# 2. Increase the left indent from 120 twips (6 pt) to 225 twips (11.25 pt).
$p1.Format.LeftIndent = 11.25

# 3. Locate the placeholder text's range dynamically via Find, scoped to
#    paragraph 1, so we don't depend on hard-coded character offsets.
$idRange = $p1.Range.Duplicate
$idRange.Find.Execute("**ID__AFFARS_pgi_5319_topic_5__ID**", $false, $false,
                       $false, $false, $false, $true, 1, $false, "", 0)

# 4. Delete the trailing run that contains only a single space character,
#    i.e. everything between the end of the placeholder text and the
#    paragraph mark.
$paraEnd = $p1.Range.End
$spaceRange = $d.Range($idRange.End, $paraEnd - 1)
$spaceRange.Delete()

# 5. Update the placeholder id text itself.
$idRange.Text = "**ID__AFFARS_AF_PGI_5319_705_2__ID**"
